# Applies the "login associated files" update to table_list.xlsx:
#  - Adds new "User"/"Account" columns (and Email/Password sub-rows) to Sheet1
#  - Updates the print/page setup for Sheet1 (A4, portrait)
#  - Moves the active selection to I3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header/labels for the "Account" block (columns G:H).
# Insert in the same order the shared strings table gains them:
#   13=User, 14=Email, 15=Password, 16=Account
$ws.Range("G1").Value = "User"
$ws.Range("H2").Value = "Email"
$ws.Range("H3").Value = "Password"
$ws.Range("H1").Value = "Account"

# Page setup: A4 paper, portrait orientation
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Move the active cell/selection to I3
$ws.Range("I3").Select() | Out-Null
